$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D, E, G columns stay formatted as text so values are written verbatim
$ws.Range("D2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Range("D2").Value = "304.01"
$ws.Range("E2").Value = "0.84%"
$ws.Range("G2").Value = "11"

$ws.Range("D3").Value = "35.69"
$ws.Range("E3").Value = "1.58%"
$ws.Range("G3").Value = "11"

$ws.Range("D4").Value = "5.082"
$ws.Range("E4").Value = "0.94%"
$ws.Range("G4").Value = "11"

$ws.Range("E5").Value = "0.87%"
$ws.Range("G5").Value = "11"

$ws.Range("D6").Value = "1.904"
$ws.Range("E6").Value = "-0.78%"
$ws.Range("G6").Value = "11"

$ws.Range("D7").Value = "7.741"
$ws.Range("E7").Value = "-0.81%"
$ws.Range("G7").Value = "11"

$ws.Range("D8").Value = "0.9281"
$ws.Range("E8").Value = "0.70%"
$ws.Range("G8").Value = "11"

$ws.Range("D9").Value = "0.1388"
$ws.Range("E9").Value = "4.98%"
$ws.Range("G9").Value = "11"

$ws.Range("D10").Value = "0.1895"
$ws.Range("E10").Value = "2.53%"
$ws.Range("G10").Value = "11"

$ws.Range("D11").Value = "0.09132"
$ws.Range("E11").Value = "-5.66%"
$ws.Range("G11").Value = "11"

$ws.Range("D12").Value = "0.03607"
$ws.Range("E12").Value = "0.61%"
$ws.Range("G12").Value = "11"

$ws.Range("D13").Value = "0.09811"
$ws.Range("E13").Value = "-0.36%"
$ws.Range("G13").Value = "11"

$ws.Range("D14").Value = "0.001406"
$ws.Range("E14").Value = "0.82%"
$ws.Range("G14").Value = "11"

$ws.Range("D15").Value = "0.005911"
$ws.Range("E15").Value = "2.52%"
$ws.Range("G15").Value = "11"

$ws.Range("D16").Value = "3.552"
$ws.Range("E16").Value = "1.22%"
$ws.Range("G16").Value = "11"

$ws.Range("D17").Value = "4.186"
$ws.Range("E17").Value = "3.52%"
$ws.Range("G17").Value = "11"

$ws.Range("D18").Value = "2.977"
$ws.Range("E18").Value = "-0.42%"
$ws.Range("G18").Value = "11"

$ws.Range("D19").Value = "0.3455"
$ws.Range("E19").Value = "0.86%"
$ws.Range("G19").Value = "11"

$ws.Range("D20").Value = "0.1332"
$ws.Range("E20").Value = "1.67%"
$ws.Range("G20").Value = "11"

$ws.Range("E21").Value = "-2.87%"
$ws.Range("G21").Value = "11"

$ws.Range("E22").Value = "2.06%"
$ws.Range("G22").Value = "11"

$ws.Range("D23").Value = "0.04450"
$ws.Range("E23").Value = "-1.03%"
$ws.Range("G23").Value = "11"

$ws.Range("D24").Value = "0.001223"
$ws.Range("E24").Value = "0.33%"
$ws.Range("G24").Value = "11"

$ws.Range("D25").Value = "0.004785"
$ws.Range("E25").Value = "-0.12%"
$ws.Range("G25").Value = "11"

$ws.Range("D26").Value = "0.0001564"
$ws.Range("E26").Value = "25.06%"
$ws.Range("G26").Value = "11"

$ws.Range("D27").Value = "0.0003135"
$ws.Range("E27").Value = "4.40%"
$ws.Range("G27").Value = "11"

$ws.Range("G28").Value = "11"

$ws.Range("G29").Value = "11"

$ws.Range("G30").Value = "11"

$ws.Range("G31").Value = "11"

$ws.Range("G32").Value = "11"

$ws.Range("G33").Value = "11"

$ws.Range("G34").Value = "11"

$ws.Range("G35").Value = "11"

$ws.Range("G36").Value = "11"

$ws.Range("G37").Value = "11"

$ws.Range("G38").Value = "11"

$ws.Range("D39").Value = "0.01961"
$ws.Range("E39").Value = "4.38%"
$ws.Range("G39").Value = "11"

$ws.Range("D40").Value = "0.04881"
$ws.Range("E40").Value = "3.80%"
$ws.Range("G40").Value = "11"

$ws.Range("D41").Value = "0.007653"
$ws.Range("E41").Value = "1.54%"
$ws.Range("G41").Value = "11"

$ws.Range("D42").Value = "0.009288"
$ws.Range("E42").Value = "-4.61%"
$ws.Range("G42").Value = "11"

$ws.Range("D43").Value = "0.1372"
$ws.Range("E43").Value = "3.77%"
$ws.Range("G43").Value = "11"

$ws.Range("E44").Value = "-0.26%"
$ws.Range("G44").Value = "11"

$ws.Range("E45").Value = "17.74%"
$ws.Range("G45").Value = "11"

$ws.Range("D46").Value = "0.00006404"
$ws.Range("E46").Value = "2.66%"
$ws.Range("G46").Value = "11"

$ws.Range("E47").Value = "0.12%"
$ws.Range("G47").Value = "11"

$ws.Range("D48").Value = "63.57"
$ws.Range("E48").Value = "-1.41%"
$ws.Range("G48").Value = "11"

$ws.Range("D49").Value = "0.001193"
$ws.Range("E49").Value = "-19.91%"
$ws.Range("G49").Value = "11"

$ws.Range("D50").Value = "0.00002104"
$ws.Range("E50").Value = "0.12%"
$ws.Range("G50").Value = "11"

$ws.Range("D51").Value = "0.0002004"
$ws.Range("E51").Value = "0.12%"
$ws.Range("G51").Value = "11"
